$wb = $excel.ActiveWorkbook

# --- Step 1: add a brand-new "Sheet3" right after "Sheet2" ---------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# --- Step 2: move the OLD Sheet2 data (Double/string/integer sample) into
#             the new Sheet3, shifted one column to the right (B->C, C->D, D->E)
$ws3.Range("C1").Value = "Double"
$ws3.Range("D1").Value = "string"
$ws3.Range("E1").Value = "integer"

$ws3.Range("C2").Value = 0.76
$ws3.Range("D2").Value = "str8"
$ws3.Range("E2").Value = 1700

$ws3.Range("C3").Value = 1.2
$ws3.Range("D3").Value = "str11"
$ws3.Range("E3").Value = 12

# --- Step 3: overwrite Sheet2 with its new content ------------------------
$ws2.Range("B1").Value = "sheet2_header1"
$ws2.Range("C1").Value = "sheet2_header2"
$ws2.Range("D1").Value = "sheet2_header3"

$ws2.Range("B2").Value = "str7"
$ws2.Range("C2").Value = "str8"
$ws2.Range("D2").Value = "str9"

$ws2.Range("B3").Value = "str10"
$ws2.Range("C3").Value = "str11"
$ws2.Range("D3").Value = "str12"

# --- Step 4: selections -----------------------------------------------
# Sheet2: select B1:D3, active cell B1, and it is no longer the tab shown
[void]$ws2.Range("B1:D3").Select()

# Sheet3: select K21 and make it the active/visible tab
[void]$ws3.Range("K21").Select()
[void]$ws3.Activate()
